$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# Insert new row 31 ("largest area in histogram") directly under row 30
# *before* row 30 itself gets touched, so the insert picks up row 30's
# original formatting (B:->s7 like B26/B29, C:->s2 wrapped style); D then
# gets its format stripped back to plain and the spurious E cell that tags
# along is cleared out completely.
# ---------------------------------------------------------------------------
$ws.Rows(31).Insert(-4121) | Out-Null                     # xlShiftDown

$ws.Range("A31").Value = "backtracking/Recursion"
$ws.Range("B31").Value = "largest area in histogram"
$ws.Range("C31").Value = "general backtracking logic/ to improve TC use hashset/(boolean - preferred easy to explain & understand) array or swapping logic"
$ws.Range("D31").Value = "https://leetcode.com/explore/learn/card/recursion-ii/507/beyond-recursion/2903/"

$ws.Range("D31").ClearFormats() | Out-Null
$ws.Range("D31").Value = "https://leetcode.com/explore/learn/card/recursion-ii/507/beyond-recursion/2903/"
$ws.Range("E31").Clear() | Out-Null

# ---------------------------------------------------------------------------
# Insert new row 32 ("Letter Combinations of a Phone Number") under row 31.
# ---------------------------------------------------------------------------
$ws.Rows(32).Insert(-4121) | Out-Null                     # xlShiftDown

$ws.Range("A32").Value = "backtracking/Recursion"
$ws.Range("B32").Value = "Letter Combinations of a Phone Number"
$ws.Range("C32").Value = "use backtracking and keypad map or an array / iterative sol uses queue FIFO"
$ws.Range("D32").Value = "https://leetcode.com/explore/learn/card/recursion-ii/507/beyond-recursion/2905/discuss/2606585/Easy-to-understand-clean-JAVA-backtracking-sol"
$ws.Range("E32").Value = "https://leetcode.com/explore/learn/card/recursion-ii/507/beyond-recursion/2905"

# B32 needs the vertical-centred style (same family as B30/D26) rather than
# the plain one it inherited from row 31.
$ws.Range("D26").Copy() | Out-Null
$ws.Range("B32").PasteSpecial(-4122) | Out-Null           # xlPasteFormats
$ws.Application.CutCopyMode = $false

# D32 carries a real hyperlink, so give it the Hyperlink cell style and wire
# up the hyperlink itself (this also creates the external relationship).
$ws.Range("D32").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("D32"), "https://leetcode.com/explore/learn/card/recursion-ii/507/beyond-recursion/2905/discuss/2606585/Easy-to-understand-clean-JAVA-backtracking-sol") | Out-Null

# ---------------------------------------------------------------------------
# Row 30: the "stack" category row gets a new B/C pairing.
#   B30 becomes "permutation" (picks up the vertical-centred font style that
#   already lives on D26 elsewhere in the sheet).
#   C30 becomes the old "trick is to find the left smaller ..." text that
#   used to live in B30/old-C30 (style stays as-is, already wrapText).
# ---------------------------------------------------------------------------
$ws.Range("C30").Value = "trick is to find the left smaller and right smaller for a bar at i that way max area for the bar is (rs -ls + 1) * bar height, to optimize compute the left & right smaller array using stack "
$ws.Range("B30").Value = "permutation"

$ws.Range("D26").Copy() | Out-Null
$ws.Range("B30").PasteSpecial(-4122) | Out-Null          # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Leave the same cell selected as the author ended up on.
$ws.Range("C32").Select() | Out-Null
